# Trade #49 (MarketMaking) closes with an early_exit at 2026-02-18 00:12:26,
# and two new trades (#78 momentum, #79 HighProbConvergence) are opened.
# This updates the Summary roll-up, the Strategy Status roll-up, the
# "All Trades" ledger, and the per-strategy ledger sheets accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: roll-up metrics after the close
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.57   # Current Capital
$summary.Range("B4").Value = 0.67      # Total P&L $
$summary.Range("B5").Value = 0.29      # Total P&L %
$summary.Range("B6").Value = 47        # Total Trades
$summary.Range("B8").Value = 18        # Losing Trades
$summary.Range("B9").Value = 55.32     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.56999999999999   # Capital
$status.Range("D6").Value = 18                  # Trades
$status.Range("E6").Value = -0.24               # P&L $
$status.Range("F6").Value = -0.43               # P&L %
$status.Range("G6").Value = 55.56               # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet: close trade #49 (row 50) + append trades #78/#79
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G50").Value = 0.681527
$allTrades.Range("H50").Value = "CLOSED"
$allTrades.Range("I50").Value = -13.7308
$allTrades.Range("J50").Value = -0.11
$allTrades.Range("K50").Value = 99.56999999999999
$allTrades.Range("L50").Value = "early_exit"
$allTrades.Range("M50").Value = 0.17

# New row 79 -> trade #78, momentum strategy, still OPEN
$row79 = @(78, "2026-02-18", "00:12:26", "momentum", "DOWN", 0.79, "", "OPEN", 0, 0, 100, "", 0, 0, 0, 0.9, "Downward momentum: -1.980% over 10 samples")
$allTrades.Cells.Item(79, 2).NumberFormat = "@"
for ($i = 0; $i -lt $row79.Length; $i++) {
    $allTrades.Cells.Item(79, $i + 1).Value = $row79[$i]
}

# New row 80 -> trade #79, HighProbConvergence strategy, still OPEN
$row80 = @(79, "2026-02-18", "00:12:26", "HighProbConvergence", "UP", 0.21, "", "OPEN", 0, 0, 100, "", 0, 0, 0, 0.95, "Mean reversion UP: price 1.79% below mean (z=-3.00)")
$allTrades.Cells.Item(80, 2).NumberFormat = "@"
for ($i = 0; $i -lt $row80.Length; $i++) {
    $allTrades.Cells.Item(80, $i + 1).Value = $row80[$i]
}

# ---------------------------------------------------------------------
# momentum sheet: append trade #78 as row 12
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentumRow = @(78, "2026-02-18", "00:12:26", "momentum", "DOWN", 0.79, "", "OPEN", 0, 0, 100, 0, 0, 0.9, "Downward momentum: -1.980% over 10 samples", "", 0)
$momentum.Cells.Item(12, 2).NumberFormat = "@"
for ($i = 0; $i -lt $momentumRow.Length; $i++) {
    $momentum.Cells.Item(12, $i + 1).Value = $momentumRow[$i]
}

# ---------------------------------------------------------------------
# HighProbConvergence sheet: append trade #79 as row 6
# ---------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpcRow = @(79, "2026-02-18", "00:12:26", "HighProbConvergence", "UP", 0.21, "", "OPEN", 0, 0, 100, 0, 0, 0.95, "Mean reversion UP: price 1.79% below mean (z=-3.00)", "", 0)
$hpc.Cells.Item(6, 2).NumberFormat = "@"
for ($i = 0; $i -lt $hpcRow.Length; $i++) {
    $hpc.Cells.Item(6, $i + 1).Value = $hpcRow[$i]
}

# ---------------------------------------------------------------------
# MarketMaking sheet: close trade #49 (row 21)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G21").Value = 0.681527
$mm.Range("H21").Value = "CLOSED"
$mm.Range("I21").Value = -13.7308
$mm.Range("J21").Value = -0.11
$mm.Range("K21").Value = 99.56999999999999
$mm.Range("P21").Value = "early_exit"
$mm.Range("Q21").Value = 0.17
